# Atualizacao dos documentos e ajustes em telas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Correct typos / wording in the requirements table (column C) ---
# RF009 - "diarios" -> "diários"
$ws.Cells.Item(11, 3).Value = "Exibição dos ganhos diários/mensais/anuais."

# RF008 - "técnina" -> "técnica"
$ws.Cells.Item(10, 3).Value = "Pré exibição de ficha de técnica."

# RF007 - add trailing period
$ws.Cells.Item(9, 3).Value = "Solicitar contraproposta do valor do serviço."

# RF013 - add trailing period
$ws.Cells.Item(15, 3).Value = "Recusar serviço."

# --- Update the view / selection state to match the saved workbook ---
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C20").Select()
